$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "FE-0317058"
$ws.Range("B5").Value = "AAA01"
$ws.Range("C5").Value = "30/07/2024 11:17:00"
$ws.Range("D5").Value = "30/07/2024 11:17:00"
$ws.Range("E5").Value = "Moto"
$ws.Range("F5").Value = 1500
$ws.Range("G5").Value = 0
$ws.Range("H5").Value = 1500
